# Fix French colon spacing in statut_name ("X: ..." -> "X : ...") and
# refresh the clinical-trials table (publipostage2/053evvt91): rows were
# re-sorted against the upstream extract, two new trials were appended
# (a EudraCT-only row for CAARDS-1, and NCT03690115 / PONALLO), and some
# acronym/title/intervention_type cells moved or were cleared accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = '3 : résultats postés ou publiés après les 36 mois'

# Row 3
$ws.Cells.Item(3, 2).Value = '3 : résultats postés ou publiés après les 36 mois'

# Row 4
$ws.Cells.Item(4, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(4, 3).Value = 'NCT02888990'
$ws.Cells.Item(4, 7).Value = 'An Open Label Phase II Study to Evaluate the Efficacy and Safety of Induction and Consolidation Therapy With Dasatinib in Combination With Chemotherapy in Patients Aged 55 Years and Over With Philadelphia Chromosome Positive (Ph+ or BCR-ABL+) Acute Lymphoblastic Leukemia (ALL).'
$ws.Cells.Item(4, 8).Value = 'EWALLPH01'
$ws.Cells.Item(4, 9).Value = 'DRUG'

# Row 5
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = '2'
$ws.Cells.Item(5, 1).Style = "Normal"
$ws.Cells.Item(5, 2).Value = '2 : résultats postés ou publiés entre 12 et 36 mois'
$ws.Cells.Item(5, 3).Value = 'NCT02883959'
$ws.Cells.Item(5, 7).ClearContents()
$ws.Cells.Item(5, 8).Value = 'Painkiller'
$ws.Cells.Item(5, 9).Value = 'OTHER'

# Row 6
$ws.Cells.Item(6, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(6, 3).Value = 'NCT02963545'
$ws.Cells.Item(6, 7).Value = 'Simultaneous Assessments of Serotonin and Kynurenine Pathways Parameters in Patients Shortly (Less Than 4 Hours and a Half) After the Onset of a Cerebral Infarction'
$ws.Cells.Item(6, 8).ClearContents()
$ws.Cells.Item(6, 9).Value = 'OTHER'

# Row 7
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = '4'
$ws.Cells.Item(7, 1).Style = "Normal"
$ws.Cells.Item(7, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(7, 3).Value = 'NCT02888964'
$ws.Cells.Item(7, 7).Value = 'A Study to Assess Efficacy and Safety of Pioglitazone as Add-On Therapy to Imatinib Mesylate in CP-CML Patients in Major Molecular Response'
$ws.Cells.Item(7, 8).Value = 'ACTIM'
$ws.Cells.Item(7, 9).Value = 'DRUG'

# Row 8
$ws.Cells.Item(8, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(8, 3).Value = 'NCT02882568'
$ws.Cells.Item(8, 7).Value = 'Analysis of Immune Response In Bacterial Infection of Obese Subject'
$ws.Cells.Item(8, 8).Value = 'ARIIBO'
$ws.Cells.Item(8, 9).Value = 'BIOLOGICAL'

# Row 9
$ws.Cells.Item(9, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(9, 3).Value = 'NCT02896842'
$ws.Cells.Item(9, 7).Value = 'A Prospective Randomized Phase II Study Evaluating the Monitoring of Imatinib Mesylate (Glivec®) Plasmatic Through Level in Patients Newly Diagnosed With Chronic Phase Chronic Myelogenous Leukaemia (CP-CML).'
$ws.Cells.Item(9, 8).Value = 'OPTIMIMATINIB'
$ws.Cells.Item(9, 9).Value = 'OTHER'

# Row 10
$ws.Cells.Item(10, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(10, 3).Value = 'NCT01946750'
$ws.Cells.Item(10, 7).ClearContents()
$ws.Cells.Item(10, 8).Value = 'SERODIFF'
$ws.Cells.Item(10, 9).Value = 'BIOLOGICAL'

# Row 11
$ws.Cells.Item(11, 2).Value = '4 : pas de résultats postés ni publiés'

# Row 12
$ws.Cells.Item(12, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(12, 3).Value = 'NCT02896374'
$ws.Cells.Item(12, 7).Value = 'Normal and Abnormal Event Related Potentials During Attribution of Intentions to Others and Resolution of Emotional Conflicts in Schizophrenia'
$ws.Cells.Item(12, 8).Value = 'SERC'
$ws.Cells.Item(12, 9).Value = 'BEHAVIORAL'

# Row 13
$ws.Cells.Item(13, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(13, 3).Value = 'NCT02901015'
$ws.Cells.Item(13, 7).Value = 'Evaluation and Validation of Social Cognition Battery to Characterize Schizophrenic Patients Functioning (EVACO)'
$ws.Cells.Item(13, 8).Value = 'EVACO'

# Row 14
$ws.Cells.Item(14, 2).Value = '2 : résultats postés ou publiés entre 12 et 36 mois'

# Row 15
$ws.Cells.Item(15, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(15, 3).Value = 'NCT03030482'
$ws.Cells.Item(15, 7).Value = 'Evaluation of Touch Massage on Anxiety in Critically Ill Patients : a Randomised Controlled Trial Study (REaLAX)'
$ws.Cells.Item(15, 8).Value = 'REaLAX'

# Row 16
$ws.Cells.Item(16, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(16, 3).Value = 'NCT02894177'
$ws.Cells.Item(16, 7).Value = 'Is Transcutaneous Carbon Dioxide Pressure (tcPCO2) Monitoring During Spontaneous Breathing Trials Useful to Predict Extubation Failure in Mechanically Ventilated Patients in the ICU?'
$ws.Cells.Item(16, 8).Value = 'tcPCO2'

# Row 17
$ws.Cells.Item(17, 2).Value = '4 : pas de résultats postés ni publiés'

# Row 18
$ws.Cells.Item(18, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(18, 3).ClearContents()
$ws.Cells.Item(18, 4).Value = '2020-001598-66'
$ws.Cells.Item(18, 6).NumberFormat = "@"
$ws.Cells.Item(18, 6).Value = '2021'
$ws.Cells.Item(18, 6).Style = "Normal"
$ws.Cells.Item(18, 7).Value = 'Preliminary randomized controlled trial of poractant alfa (Curosurf®) by fiberoptic bronchoscopy-directed endobronchial administration in acute respiratory distress syndrome (ARDS) due to COVID-19 viral pneumonia. 
 Essai thérapeutique randomisé contrôlé préliminaire du poractant alfa (Curosurf®) en administration endobronchique dirigée par fibroscopie dans le syndrome de détresse respiratoire aiguë provoqué par la pneumonie virale COVID-19'
$ws.Cells.Item(18, 8).Value = 'CAARDS-1'
$ws.Cells.Item(18, 9).Value = 'DRUG (presumed)'

# Row 19
$ws.Cells.Item(19, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(19, 3).Value = 'NCT03973385'
$ws.Cells.Item(19, 6).NumberFormat = "@"
$ws.Cells.Item(19, 6).Value = '2022'
$ws.Cells.Item(19, 6).Style = "Normal"
$ws.Cells.Item(19, 7).Value = 'Evaluation of Efficacy of Cryotherapy for Skin Anesthesia During Arterial Puncture for Blood Gas in Critically Ill Patients. SNOW Study'
$ws.Cells.Item(19, 8).Value = 'SNOW'
$ws.Cells.Item(19, 9).Value = 'DEVICE'

# Row 20
$ws.Cells.Item(20, 1).NumberFormat = "@"
$ws.Cells.Item(20, 1).Value = '4'
$ws.Cells.Item(20, 1).Style = "Normal"
$ws.Cells.Item(20, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(20, 3).Value = 'NCT02473146'
$ws.Cells.Item(20, 6).NumberFormat = "@"
$ws.Cells.Item(20, 6).Value = '2023'
$ws.Cells.Item(20, 6).Style = "Normal"
$ws.Cells.Item(20, 7).Value = 'Etude Exploratoire randomisée Comparant le Traitement Par Gemtuzumab Ozogamicin /Cytarabine au Traitement Standard Par Idarubicine/Cytarabinechez Les Sujets âgés de 60 à 80 Ans et présentant Une LAM et un Caryotype Non défavorable'
$ws.Cells.Item(20, 8).Value = 'ALFA1401'
$ws.Cells.Item(20, 9).Value = 'DRUG'

# Row 21
$ws.Cells.Item(21, 1).NumberFormat = "@"
$ws.Cells.Item(21, 1).Value = '4'
$ws.Cells.Item(21, 1).Style = "Normal"
$ws.Cells.Item(21, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(21, 3).Value = 'NCT03690115'
$ws.Cells.Item(21, 6).NumberFormat = "@"
$ws.Cells.Item(21, 6).Value = '2024'
$ws.Cells.Item(21, 6).Style = "Normal"
$ws.Cells.Item(21, 7).Value = 'Phase 2 Study of Ponatinib (Iclusig) for Prevention of Relapse After Allogeneic Stem Cell Transplantation (allo-SCT) in FLT3-ITD AML Patients: the PONALLO Trial."'
$ws.Cells.Item(21, 8).Value = 'PONALLO'
$ws.Cells.Item(21, 9).Value = 'DRUG'
